$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @{
    "B2" = 0.1764705882352941;
    "C2" = 0.6470588235294118;
    "P2" = 0.1764705882352941;
    "J3" = 0.09090909090909091;
    "P3" = 0.7272727272727273;
    "S3" = 0.1818181818181818;
    "P4" = 0.25;
    "S4" = 0.75;
    "J6" = 0.2142857142857143;
    "Q6" = 0.4285714285714285;
    "S6" = 0.3571428571428572;
    "B7" = 0.1428571428571428;
    "F7" = 0.1428571428571428;
    "O7" = 0.1428571428571428;
    "Q7" = 0.2857142857142857;
    "S7" = 0.2857142857142857;
    "B8" = 0.1052631578947368;
    "D8" = 0.05263157894736842;
    "F8" = 0.1052631578947368;
    "J8" = 0.2105263157894737;
    "Q8" = 0.2105263157894737;
    "R8" = 0.1052631578947368;
    "S8" = 0.2105263157894737;
    "B9" = 0.04545454545454546;
    "F9" = 0.04545454545454546;
    "J9" = 0.1818181818181818;
    "Q9" = 0.3636363636363636;
    "R9" = 0.1363636363636364;
    "S9" = 0.2272727272727273;
    "B10" = 0.08928571428571429;
    "D10" = 0.01785714285714286;
    "F10" = 0.04464285714285714;
    "J10" = 0.1160714285714286;
    "O10" = 0.01785714285714286;
    "Q10" = 0.3482142857142857;
    "R10" = 0.07142857142857142;
    "S10" = 0.2946428571428572;
    "G11" = 0.1333333333333333;
    "J11" = 0.2;
    "K11" = 0.2;
    "L11" = 0.4666666666666667;
    "G12" = 0.5;
    "J12" = 0.375;
    "L12" = 0.125;
    "G13" = 1;
    "I15" = 0.1333333333333333;
    "J15" = 0.2;
    "K15" = 0.1333333333333333;
    "S15" = 0.5333333333333333;
    "J16" = 0.6363636363636364;
    "O16" = 0.09090909090909091;
    "S16" = 0.2727272727272727;
    "H17" = 0.1186440677966102;
    "I17" = 0.2033898305084746;
    "J17" = 0.4745762711864407;
    "K17" = 0.05084745762711865;
    "M17" = 0.01694915254237288;
    "O17" = 0.03389830508474576;
    "S17" = 0.1016949152542373;
    "J18" = 0.75;
    "K18" = 0.08333333333333333;
    "O18" = 0.08333333333333333;
    "S18" = 0.08333333333333333;
    "F19" = 0.02439024390243903;
    "H19" = 0.1463414634146341;
    "I19" = 0.0975609756097561;
    "J19" = 0.4268292682926829;
    "K19" = 0.07317073170731707;
    "O19" = 0.08536585365853659;
    "S19" = 0.1463414634146341
}

foreach ($cellRef in $values.Keys) {
    $ws.Range($cellRef).Value = $values[$cellRef]
}
